# Auto-generated edit script applying the cell-value changes described by the diff.
# Each sheet's changes are applied via direct Range(...).Value assignments;
# cells that are removed entirely in the target state use ClearContents().
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 789.2857
$ws.Range("I9").Value = 155
$ws.Range("J9").Value = 2375
$ws.Range("K9").Value = 155
$ws.Range("L9").Value = 2375
$ws.Range("M9").Value = 14
$ws.Range("N9").Value = -2713
$ws.Range("H40").Value = 3233
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3233
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 3233
$ws.Range("N40").Value = -3583
$ws.Range("H43").Value = 5007
$ws.Range("I43").Value = 3666.6667
$ws.Range("J43").Value = 6155.857
$ws.Range("K43").Value = 3666.6667
$ws.Range("L43").Value = 6155.857
$ws.Range("M43").Value = -3597.6667
$ws.Range("N43").Value = -6293.857
$ws.Range("H53").Value = 537.5
$ws.Range("I53").Value = 134.66667
$ws.Range("K53").Value = 134.66667
$ws.Range("M53").Value = 502.33333
$ws.Range("H99").Value = 1161.25
$ws.Range("I99").Value = 732
$ws.Range("J99").Value = 1590.5
$ws.Range("K99").Value = 2196
$ws.Range("L99").Value = 4771.5
$ws.Range("M99").Value = -698
$ws.Range("N99").Value = -7767.5
$ws.Range("H103").Value = 2785.158
$ws.Range("J103").Value = 4811.8
$ws.Range("L103").Value = 14435.4
$ws.Range("N103").Value = -15607.4
$ws.Range("H135").Value = 2469.7144
$ws.Range("J135").Value = 4996.3335
$ws.Range("L135").Value = 44967.0015
$ws.Range("N135").Value = -50037.0015
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("N136").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 3501
$ws.Range("I23").Value = 3501
$ws.Range("K23").Value = 3501
$ws.Range("M23").Value = -3242
$ws.Range("H88").Value = 3687.5
$ws.Range("I88").Value = 4502.5
$ws.Range("J88").Value = 3415.8333
$ws.Range("K88").Value = 4502.5
$ws.Range("L88").Value = 3415.8333
$ws.Range("M88").Value = -4096.5
$ws.Range("N88").Value = -4227.8333
$ws.Range("H91").Value = 3687.5
$ws.Range("I91").Value = 4502.5
$ws.Range("J91").Value = 3415.8333
$ws.Range("K91").Value = 4502.5
$ws.Range("L91").Value = 3415.8333
$ws.Range("M91").Value = -3098.5
$ws.Range("N91").Value = -6223.8333
$ws.Range("H134").Value = 79974
$ws.Range("J134").Value = 79974
$ws.Range("L134").Value = 79974
$ws.Range("N134").Value = -90114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").ClearContents()
$ws.Range("N6").Value = 0
$ws.Range("H55").Value = 75233
$ws.Range("J55").Value = 75233
$ws.Range("L55").Value = 75233
$ws.Range("N55").Value = -75779
$ws.Range("H86").Value = 9199
$ws.Range("I86").Value = 7665.3335
$ws.Range("J86").Value = 11499.5
$ws.Range("K86").Value = 7665.3335
$ws.Range("L86").Value = 11499.5
$ws.Range("M86").Value = -6542.3335
$ws.Range("N86").Value = -13745.5
$ws.Range("H89").Value = 9199
$ws.Range("I89").Value = 7665.3335
$ws.Range("J89").Value = 11499.5
$ws.Range("K89").Value = 38326.6675
$ws.Range("L89").Value = 57497.5
$ws.Range("M89").Value = -32710.6675
$ws.Range("N89").Value = -68729.5
$ws.Range("H94").Value = 895.9091
$ws.Range("I94").Value = 485.5
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 485.5
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -34.5
$ws.Range("N94").Value = -5902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2689
$ws.Range("J31").Value = 3505.25
$ws.Range("L31").Value = 3505.25
$ws.Range("N31").Value = -4095.25
$ws.Range("H34").Value = 2689
$ws.Range("J34").Value = 3505.25
$ws.Range("L34").Value = 3505.25
$ws.Range("N34").Value = -3909.25
$ws.Range("H86").Value = 1524974.8
$ws.Range("I86").Value = 1599966.6
$ws.Range("J86").Value = 1299999
$ws.Range("K86").Value = 1599966.6
$ws.Range("L86").Value = 1299999
$ws.Range("M86").Value = -1598843.6
$ws.Range("N86").Value = -1302245
$ws.Range("H89").Value = 1524974.8
$ws.Range("I89").Value = 1599966.6
$ws.Range("J89").Value = 1299999
$ws.Range("K89").Value = 7999833
$ws.Range("L89").Value = 6499995
$ws.Range("M89").Value = -7994217
$ws.Range("N89").Value = -6511227
$ws.Range("H105").Value = 6446.778
$ws.Range("I105").Value = 7824.4287
$ws.Range("J105").Value = 1625
$ws.Range("K105").Value = 7824.4287
$ws.Range("L105").Value = 1625
$ws.Range("M105").Value = -6077.4287
$ws.Range("N105").Value = -5119
$ws.Range("H107").Value = 388.7647
$ws.Range("I107").Value = 364.75
$ws.Range("J107").Value = 446.4
$ws.Range("K107").Value = 364.75
$ws.Range("L107").Value = 446.4
$ws.Range("M107").Value = 1555.25
$ws.Range("N107").Value = -4286.4
$ws.Range("H115").Value = 45000
$ws.Range("I115").Value = 45000
$ws.Range("K115").Value = 45000
$ws.Range("M115").Value = -43825
$ws.Range("H121").Value = 37000
$ws.Range("J121").Value = 37000
$ws.Range("L121").Value = 37000
$ws.Range("N121").Value = -39620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 931.6667
$ws.Range("I5").Value = 901
$ws.Range("J5").Value = 947
$ws.Range("K5").Value = 2703
$ws.Range("L5").Value = 2841
$ws.Range("M5").Value = -2591
$ws.Range("N5").Value = -3065
$ws.Range("H14").Value = 2499.5
$ws.Range("I14").Value = 2499.5
$ws.Range("K14").Value = 7498.5
$ws.Range("M14").Value = -7325.5
$ws.Range("H41").Value = 140
$ws.Range("I41").Value = 140
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 420
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -82
$ws.Range("H93").Value = 2707.5
$ws.Range("I93").Value = 249
$ws.Range("J93").Value = 15000
$ws.Range("K93").Value = 747
$ws.Range("L93").Value = 45000
$ws.Range("M93").Value = 1125
$ws.Range("N93").Value = -48744
$ws.Range("H131").Value = 894.9524
$ws.Range("I131").Value = 959.3333
$ws.Range("J131").Value = 890
$ws.Range("K131").Value = 2877.9999
$ws.Range("L131").Value = 2670
$ws.Range("M131").Value = 2162.0001
$ws.Range("N131").Value = -12750
$ws.Range("H135").Value = 931.6667
$ws.Range("I135").Value = 901
$ws.Range("J135").Value = 947
$ws.Range("K135").Value = 8109
$ws.Range("L135").Value = 8523
$ws.Range("M135").Value = -5574
$ws.Range("N135").Value = -13593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12000
$ws.Range("I70").Value = 12000
$ws.Range("K70").Value = 12000
$ws.Range("M70").Value = -11730
$ws.Range("H73").Value = 12000
$ws.Range("I73").Value = 12000
$ws.Range("K73").Value = 12000
$ws.Range("M73").Value = -11064
$ws.Range("H102").Value = 2455.0908
$ws.Range("I102").Value = 2289.5
$ws.Range("K102").Value = 2289.5
$ws.Range("M102").Value = -667.5
$ws.Range("H122").Value = 2124.75
$ws.Range("I122").Value = 2124.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6374.25
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -3924.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1079.5714
$ws.Range("I35").Value = 1079.5714
$ws.Range("K35").Value = 1079.5714
$ws.Range("M35").Value = -743.5714
$ws.Range("I46").Value = 4111.1113
$ws.Range("J46").Value = 4677.4194
$ws.Range("K46").Value = 4111.1113
$ws.Range("L46").Value = 4677.4194
$ws.Range("M46").Value = -3923.1113
$ws.Range("N46").Value = -5053.4194
$ws.Range("H55").Value = 469.73334
$ws.Range("J55").Value = 535.0769
$ws.Range("L55").Value = 535.0769
$ws.Range("N55").Value = -881.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0
$ws.Range("H6").Value = 766.3333
$ws.Range("I6").Value = 649.5
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 649.5
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = -534.5
$ws.Range("N6").Value = -1230
$ws.Range("H118").Value = 24166.666
$ws.Range("I118").Value = 20000
$ws.Range("J118").Value = 25000
$ws.Range("K118").Value = 20000
$ws.Range("L118").Value = 25000
$ws.Range("M118").Value = -18343
$ws.Range("N118").Value = -28314
